# Arbeitszeiten.xlsx edit: "fixed encoding + little frontend bug"
# - Adds a David row ("Testing + fixing", 27.11.2018, 2h) and a David row
#   ("Protokoll", 01.12.2018, 1h) before the existing 29.11.2018 Daniel entry.
# - Appends a final David row ("Encoding fixes", 14.12.2018, 5h) at the end.
#
# Rows are shifted manually (value-by-value) instead of via Rows.Insert()
# so that the existing SUMIF formulas in G3/H3 keep referencing B2:B149 /
# D2:D149 exactly as before (only their cached results change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateStyleSource = $ws.Range("A2")

function Set-DateCell($cell, $serial) {
    $cell.Value = $serial
    $dateStyleSource.Copy()
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Pre-seed the three new shared strings in the exact order they need
#    to appear in sharedStrings.xml: "Encoding fixes", "Testing + fixing",
#    "Protokoll". Using far-away scratch cells (outside any used range)
#    guarantees the row-shift below (step 2) does not move them before
#    they are cleared again.
# ---------------------------------------------------------------------
$scratch1 = $ws.Cells.Item(500, 1)
$scratch2 = $ws.Cells.Item(501, 1)
$scratch3 = $ws.Cells.Item(502, 1)

$scratch1.Value = "Encoding fixes"
$scratch2.Value = "Testing + fixing"
$scratch3.Value = "Protokoll"

$scratch1.Value = ""
$scratch2.Value = ""
$scratch3.Value = ""

# ---------------------------------------------------------------------
# 2) Shift existing rows 65-72 down to rows 67-74 (bottom-up, so no data
#    is overwritten before it is read).
# ---------------------------------------------------------------------
for ($r = 72; $r -ge 65; $r--) {
    $dest = $r + 2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2()
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value2()
    Set-DateCell $ws.Cells.Item($dest, 1) $ws.Cells.Item($r, 1).Value2()
}

# ---------------------------------------------------------------------
# 3) Fill in the two newly freed rows (65, 66) with the new entries.
# ---------------------------------------------------------------------
Set-DateCell $ws.Cells.Item(65, 1) 43431
$ws.Cells.Item(65, 2).Value = "David"
$ws.Cells.Item(65, 3).Value = "Testing + fixing"
$ws.Cells.Item(65, 4).Value = 2

Set-DateCell $ws.Cells.Item(66, 1) 43435
$ws.Cells.Item(66, 2).Value = "David"
$ws.Cells.Item(66, 3).Value = "Protokoll"
$ws.Cells.Item(66, 4).Value = 1

# ---------------------------------------------------------------------
# 4) Append the new final row (75) after the shifted data (now ending at
#    row 74).
# ---------------------------------------------------------------------
Set-DateCell $ws.Cells.Item(75, 1) 43448
$ws.Cells.Item(75, 2).Value = "David"
$ws.Cells.Item(75, 3).Value = "Encoding fixes"
$ws.Cells.Item(75, 4).Value = 5

# ---------------------------------------------------------------------
# 5) Refresh the view state: active cell / scroll position.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C76").Select()
